# Adds 6 new "CR0_*" / "XRA_*" indicator rows to the
# "r AnalysisUnit_Variable" sheet, mirroring the existing rows 3-52.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("r AnalysisUnit_Variable")

$newVariables = @(
    "CR0_TUA",
    "CR0_TUR",
    "CR0_TUS",
    "CR0_UCFBT",
    "XRA_GSC_MAX_L3M",
    "XRA_RSC_ACC_L1M "
)

$firstNewRow = 53

for ($i = 0; $i -lt $newVariables.Length; $i++) {
    $rowNum = $firstNewRow + $i

    # Copy the row directly above (which carries the correct look & feel /
    # style of the existing data rows) and insert it as a new row, shifting
    # everything below it down.
    $ws.Rows.Item($rowNum - 1).Copy()
    $ws.Rows.Item($rowNum).Insert(-4121)

    $name = $newVariables[$i]

    $ws.Cells.Item($rowNum, 1).Value2 = "CREATE/MODIFY"
    $ws.Cells.Item($rowNum, 2).Value2 = $name
    $ws.Cells.Item($rowNum, 3).Value2 = $name
    $ws.Cells.Item($rowNum, 5).Value2 = "CUSTOMER"
    $ws.Cells.Item($rowNum, 6).Value2 = $name
}

# Match the saved view state: scrolled near the bottom of the list with the
# last populated "Variable" cell selected.
[void]$ws.Activate()
[void]$ws.Range("F55").Select()
